# FOML_215_PPT.pptx edit script
# 1) Bump the cached "today" date shown on the Notes Master (Insert > Header &
#    Footer date placeholder) from 08-05-2025 to 09-05-2025.
# 2) Remove the leftover "Second/Zeroth Review" Date placeholder that had been
#    added on top of every content slide (slides 2-14) - done the same way a
#    user would in PowerPoint, via Insert > Header & Footer > uncheck Date,
#    which drops the placeholder shape from each slide.
# 3) Re-title slide 8 from "List of Modules" to "List of Frameworks".

$p = $ppt.ActivePresentation

# --- 1) Notes Master date field -------------------------------------------
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.Name -eq "Date Placeholder 2") {
        $shp.TextFrame.TextRange.Text = "09-05-2025"
    }
}

# --- 2) Drop the per-slide Date placeholder on slides 2 through 14 --------
for ($i = 2; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $slide.HeadersFooters.DateAndTime.Visible = $false
}

# --- 3) Update slide 8's title text ----------------------------------------
$slide8 = $p.Slides.Item(8)
for ($i = 1; $i -le $slide8.Shapes.Count; $i++) {
    $shp = $slide8.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "List of Modules") {
        $shp.TextFrame.TextRange.Text = "List of Frameworks"
    }
}
